$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 4981.8335
$ws.Range("J32").Value = 7196.3335
$ws.Range("L32").Value = 7196.3335
$ws.Range("N32").Value = -7848.3335
$ws.Range("H33").Value = 435.3125
$ws.Range("I33").Value = 447.66666
$ws.Range("J33").Value = 250
$ws.Range("K33").Value = 447.66666
$ws.Range("L33").Value = 250
$ws.Range("M33").Value = -218.66666
$ws.Range("N33").Value = -708
$ws.Range("H64").Value = 4624.5
$ws.Range("I64").Value = 4999
$ws.Range("J64").Value = 4250
$ws.Range("K64").Value = 4999
$ws.Range("L64").Value = 4250
$ws.Range("M64").Value = -4751
$ws.Range("N64").Value = -4746
$ws.Range("H67").Value = 4624.5
$ws.Range("I67").Value = 4999
$ws.Range("J67").Value = 4250
$ws.Range("K67").Value = 4999
$ws.Range("L67").Value = 4250
$ws.Range("M67").Value = -4141
$ws.Range("N67").Value = -5966
$ws.Range("H106").Value = 4480.4
$ws.Range("I106").Value = 4743.357
$ws.Range("J106").Value = 799
$ws.Range("K106").Value = 4743.357
$ws.Range("L106").Value = 799
$ws.Range("M106").Value = -4112.357
$ws.Range("N106").Value = -2061

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2306.68
$ws.Range("I32").Value = 2375.1667
$ws.Range("K32").Value = 2375.1667
$ws.Range("M32").Value = -2088.1667
$ws.Range("H61").Value = 5251.763
$ws.Range("I61").Value = 4472.8066
$ws.Range("J61").Value = 8701.429
$ws.Range("K61").Value = 4472.8066
$ws.Range("L61").Value = 8701.429
$ws.Range("M61").Value = -4260.8066
$ws.Range("N61").Value = -9125.429
$ws.Range("H74").Value = 2057.9092
$ws.Range("I74").Value = 2057.9092
$ws.Range("K74").Value = 2057.9092
$ws.Range("M74").Value = -1183.9092
$ws.Range("H77").Value = 2057.9092
$ws.Range("I77").Value = 2057.9092
$ws.Range("K77").Value = 10289.546
$ws.Range("M77").Value = -5921.546
$ws.Range("H122").Value = 1364.7858
$ws.Range("I122").Value = 1394
$ws.Range("K122").Value = 4182
$ws.Range("M122").Value = -1732
$ws.Range("H132").Value = 2294.4285
$ws.Range("I132").Value = 2294.4285
$ws.Range("K132").Value = 6883.2855
$ws.Range("M132").Value = -4353.2855
$ws.Range("H136").Value = 5251.763
$ws.Range("I136").Value = 4472.8066
$ws.Range("J136").Value = 8701.429
$ws.Range("K136").Value = 13418.4198
$ws.Range("L136").Value = 26104.287
$ws.Range("M136").Value = -10868.4198
$ws.Range("N136").Value = -31204.287

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2729.6667
$ws.Range("I20").Value = 2513.8
$ws.Range("K20").Value = 2513.8
$ws.Range("M20").Value = -2266.8
$ws.Range("H134").Value = 3472.4167
$ws.Range("I134").Value = 3317.0571
$ws.Range("K134").Value = 9951.1713
$ws.Range("M134").Value = -7416.1713

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 62489.5
$ws.Range("J20").Value = 62489.5
$ws.Range("L20").Value = 62489.5
$ws.Range("N20").Value = -62961.5
$ws.Range("H30").Value = 62489.5
$ws.Range("J30").Value = 62489.5
$ws.Range("L30").Value = 62489.5
$ws.Range("N30").Value = -62671.5
$ws.Range("H31").Value = 2368.4375
$ws.Range("I31").Value = 2117.9092
$ws.Range("J31").Value = 2919.6
$ws.Range("K31").Value = 2117.9092
$ws.Range("L31").Value = 2919.6
$ws.Range("M31").Value = -1822.9092
$ws.Range("N31").Value = -3509.6
$ws.Range("H34").Value = 2368.4375
$ws.Range("I34").Value = 2117.9092
$ws.Range("J34").Value = 2919.6
$ws.Range("K34").Value = 2117.9092
$ws.Range("L34").Value = 2919.6
$ws.Range("M34").Value = -1915.9092
$ws.Range("N34").Value = -3323.6
$ws.Range("H62").Value = 8499.666999999999
$ws.Range("I62").Value = 8500
$ws.Range("J62").Value = 8499.5
$ws.Range("K62").Value = 8500
$ws.Range("L62").Value = 8499.5
$ws.Range("M62").Value = -7876
$ws.Range("N62").Value = -9747.5
$ws.Range("H65").Value = 8499.666999999999
$ws.Range("I65").Value = 8500
$ws.Range("J65").Value = 8499.5
$ws.Range("K65").Value = 42500
$ws.Range("L65").Value = 42497.5
$ws.Range("M65").Value = -39380
$ws.Range("N65").Value = -48737.5
$ws.Range("H107").Value = 668.5714
$ws.Range("I107").Value = 668.5714
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 668.5714
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1251.4286
$ws.Range("N107").ClearContents()
$ws.Range("H128").Value = 62489.5
$ws.Range("J128").Value = 62489.5
$ws.Range("L128").Value = 62489.5
$ws.Range("N128").Value = -72449.5
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()
$ws.Range("H134").Value = 1512.0233
$ws.Range("I134").Value = 1412.95
$ws.Range("K134").Value = 4238.85
$ws.Range("M134").Value = -1703.85

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H49").Value = 1186.25
$ws.Range("I49").Value = 1272.5
$ws.Range("J49").Value = 1100
$ws.Range("K49").Value = 3817.5
$ws.Range("L49").Value = 3300
$ws.Range("M49").Value = -3661.5
$ws.Range("N49").Value = -3612
$ws.Range("H74").Value = 750
$ws.Range("I74").Value = 750
$ws.Range("K74").Value = 2250
$ws.Range("M74").Value = -1189
$ws.Range("H77").Value = 750
$ws.Range("I77").Value = 750
$ws.Range("K77").Value = 6750
$ws.Range("M77").Value = -1446
$ws.Range("H86").Value = 1697.5
$ws.Range("I86").Value = 764.5
$ws.Range("K86").Value = 2293.5
$ws.Range("M86").Value = -1107.5
$ws.Range("H89").Value = 1697.5
$ws.Range("I89").Value = 764.5
$ws.Range("K89").Value = 6880.5
$ws.Range("M89").Value = -952.5
$ws.Range("H97").Value = 955.7143
$ws.Range("J97").Value = 997
$ws.Range("L97").Value = 2991
$ws.Range("N97").Value = -3983
$ws.Range("H132").Value = 2332.6667
$ws.Range("I132").Value = 1999
$ws.Range("J132").Value = 2499.5
$ws.Range("K132").Value = 17991
$ws.Range("L132").Value = 22495.5
$ws.Range("M132").Value = -15461
$ws.Range("N132").Value = -27555.5

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3312.25
$ws.Range("I102").Value = 3312.25
$ws.Range("K102").Value = 3312.25
$ws.Range("M102").Value = -1690.25
$ws.Range("H122").Value = 3357.0588
$ws.Range("I122").Value = 3656.1667
$ws.Range("K122").Value = 10968.5001
$ws.Range("M122").Value = -8518.500100000001
$ws.Range("H126").Value = 4079.5
$ws.Range("I126").Value = 2772
$ws.Range("K126").Value = 8316
$ws.Range("M126").Value = -5846
$ws.Range("H131").Value = 43774.668
$ws.Range("J131").Value = 43774.668
$ws.Range("L131").Value = 43774.668
$ws.Range("N131").Value = -53854.668
$ws.Range("H132").Value = 1646.1154
$ws.Range("I132").Value = 1591.96
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 4775.88
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -2245.88
$ws.Range("N132").Value = -14060

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5752.9165
$ws.Range("I7").Value = 5621.364
$ws.Range("K7").Value = 5621.364
$ws.Range("M7").Value = -5509.364
$ws.Range("H36").Value = 70000
$ws.Range("J36").Value = 70000
$ws.Range("L36").Value = 70000
$ws.Range("N36").Value = -71124
$ws.Range("H40").Value = 1649.75
$ws.Range("I40").Value = 1649.75
$ws.Range("K40").Value = 1649.75
$ws.Range("M40").Value = -1513.75
$ws.Range("H68").Value = 2582.9167
$ws.Range("I68").Value = 2412
$ws.Range("J68").Value = 2924.75
$ws.Range("K68").Value = 2412
$ws.Range("L68").Value = 2924.75
$ws.Range("M68").Value = -1663
$ws.Range("N68").Value = -4422.75
$ws.Range("H71").Value = 2582.9167
$ws.Range("I71").Value = 2412
$ws.Range("J71").Value = 2924.75
$ws.Range("K71").Value = 12060
$ws.Range("L71").Value = 14623.75
$ws.Range("M71").Value = -8316
$ws.Range("N71").Value = -22111.75
$ws.Range("H126").Value = 5752.9165
$ws.Range("I126").Value = 5621.364
$ws.Range("K126").Value = 16864.092
$ws.Range("M126").Value = -14394.092
$ws.Range("H132").Value = 2712.158
$ws.Range("I132").Value = 2696
$ws.Range("J132").Value = 2849.5
$ws.Range("K132").Value = 8088
$ws.Range("L132").Value = 8548.5
$ws.Range("M132").Value = -5558
$ws.Range("N132").Value = -13608.5
$ws.Range("H136").Value = 4388.087
$ws.Range("I136").Value = 3134.889
$ws.Range("K136").Value = 9404.667000000001
$ws.Range("M136").Value = -6854.667000000001

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 45172
$ws.Range("I45").Value = 37872.145
$ws.Range("K45").Value = 37872.145
$ws.Range("M45").Value = -37381.145
$ws.Range("H81").Value = 3239.4614
$ws.Range("I81").Value = 3857.2222
$ws.Range("J81").Value = 1849.5
$ws.Range("K81").Value = 7714.4444
$ws.Range("L81").Value = 3699
$ws.Range("M81").Value = -6653.4444
$ws.Range("N81").Value = -5821
$ws.Range("H84").Value = 3239.4614
$ws.Range("I84").Value = 3857.2222
$ws.Range("J84").Value = 1849.5
$ws.Range("K84").Value = 38572.222
$ws.Range("L84").Value = 18495
$ws.Range("M84").Value = -33268.222
$ws.Range("N84").Value = -29103
$ws.Range("H122").Value = 6040.4443
$ws.Range("I122").Value = 5648.6
$ws.Range("J122").Value = 7999.6665
$ws.Range("K122").Value = 16945.8
$ws.Range("L122").Value = 23998.9995
$ws.Range("M122").Value = -14495.8
$ws.Range("N122").Value = -28898.9995
$ws.Range("H130").Value = 54994
$ws.Range("J130").Value = 54994
$ws.Range("L130").Value = 54994
$ws.Range("N130").Value = -65034
$ws.Range("H136").Value = 13522.692
$ws.Range("I136").Value = 11840.8
$ws.Range("J136").Value = 19129
$ws.Range("K136").Value = 35522.39999999999
$ws.Range("L136").Value = 57387
$ws.Range("M136").Value = -32972.39999999999
$ws.Range("N136").Value = -62487
